# Applies the "branch alternate" annotation edit to the first paragraph and
# appends a new shaded (empty) paragraph at the end of the document.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gets two
#    trailing spaces appended, followed by a new red (C00000) annotation
#    split across three runs:
#       "(This is a change " + EN DASH + " Ve"
#       "rsion for branch alternate"
#       ")"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$base = $d.Range($p1.Start, $p1.End - 1)   # paragraph text, excludes the paragraph mark
$base.InsertAfter("  ")

$enDash = [string][char]0x2013
$openParen = [string][char]0x0028
$closeParen = [string][char]0x0029

$run1 = $d.Range($base.End, $base.End)
$run1.InsertAfter($openParen + "This is a change " + $enDash + " Ve")
$run1.Font.Color = 192          # 0xC00000 (RGB) stored as BGR long -> 192

$run2 = $d.Range($run1.End, $run1.End)
$run2.InsertAfter("rsion for branch alternate")
$run2.Font.Color = 192

$run3 = $d.Range($run2.End, $run2.End)
$run3.InsertAfter($closeParen)
$run3.Font.Color = 192

# ---------------------------------------------------------------------
# 2) Append a brand-new, plain paragraph (Normal style, no direct run
#    formatting) shaded with fill color F9F9F9, at the very end of the
#    document body. Built via a raw OOXML fragment so it does not pick
#    up any inherited direct formatting from the preceding paragraph.
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)

$flatOpcXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr>' +
                '<w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/>' +
              '</w:pPr>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$endRange.InsertXML($flatOpcXml) | Out-Null
